$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: NewLoanInput
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# "product" value changes from "Chaithanya 123" -> "chaithanyatest"
$ws1.Range("B2").Value = "chaithanyatest"

# A new "Firstrepaymenton" row is inserted before the current row 7
# (maximumallowedoutstandingbalance), pushing everything below down by one.
$ws1.Rows.Item(7).Insert()

# Give the new row 7 the same look as the other date-label rows (e.g. the
# "expecteddisbursementon1" row, now at row 13): wrap-text label cell in A,
# date-formatted value cell in B.
$ws1.Range("A13").Copy()
$ws1.Range("A7").PasteSpecial(-4122) # xlPasteFormats
$ws1.Range("B13").Copy()
$ws1.Range("B7").PasteSpecial(-4122) # xlPasteFormats

$ws1.Range("A7").Value = "Firstrepaymenton"
$ws1.Range("B7").Value = 42036

# ---------------------------------------------------------------------------
# Sheet 4: Transactions  (Entry ID value changes from 2827 -> 193)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = 193

# ---------------------------------------------------------------------------
# Selections - restore the per-sheet active cell that the author left behind.
# Re-selecting every sheet (ending on the sheet that is really active, sheet
# 4) keeps the saved "tabSelected" flag on the correct tab.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Range("B12").Select()
$ws2.Range("D4").Select()
$ws3.Range("G5").Select()
$ws4.Range("C2").Select()

Write-Host "edit complete"
